$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.153.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.93%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.98"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4667"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2835"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06551"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.18"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -7.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.855.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.127"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6700"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.95"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.175.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.451"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.098.81"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007257"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.144"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.311"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.92"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -9.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.345"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09606"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.405"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.470"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.118"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04664"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7010"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.098"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01857"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.527"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.44"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8518"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.927"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4161"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.15"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "987.86"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.140"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.170"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1139"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.80%  "
